# Weekly fruit/hortaliza update: insert a new, more-recent price record
# (as row 28) above the existing "Poroto verde" history, pushing the
# former rows 28-47 down to 29-48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 28 (shifts 28..47 -> 29..48).
$ws.Rows.Item(28).Insert()

# Populate the new row 28 with the latest market record.
$ws.Cells.Item(28, 1).Value = 1
$ws.Cells.Item(28, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(28, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(28, 4).Value = 44567
$ws.Cells.Item(28, 5).Value = 15
$ws.Cells.Item(28, 6).Value = 100112031
$ws.Cells.Item(28, 7).Value = "Poroto verde"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 1200
$ws.Cells.Item(28, 11).Value = 400
$ws.Cells.Item(28, 12).Value = 500
$ws.Cells.Item(28, 13).Value = 450
$ws.Cells.Item(28, 14).Value = "$/kilo"
$ws.Cells.Item(28, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(28, 16).Value = 450
$ws.Cells.Item(28, 17).Value = 1
$ws.Cells.Item(28, 18).Value = "Hortaliza"
